# Apply crypto price-tracker refresh (GitHub Actions scheduled update).
# Updates Price (D) / Volume(1h) (E) columns for each row, plus a few rows
# whose ranking (and therefore Coin/Link) shifted position.
# Numeric-looking Price strings are prefixed with a literal apostrophe so
# Excel keeps them as text (preserving formats like '1.000' / '0.08586')
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '24.595.59'
$ws.Range("E2").Value = '  +2.06%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.694.97'
$ws.Range("E3").Value = '  +2.10%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.83%  '

# Row 5: BNB
$ws.Range("D5").Value = '''313.32'
$ws.Range("E5").Value = '  +1.36%  '

# Row 6: USDC
$ws.Range("D6").Value = '''1.0000'
$ws.Range("E6").Value = '  -0.62%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.3951'
$ws.Range("E7").Value = '  +1.49%  '

# Row 8: Cardano
$ws.Range("D8").Value = '''0.4034'
$ws.Range("E8").Value = '  +2.62%  '

# Row 9: BinanceUSD->OKB
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '''56.77'
$ws.Range("E9").Value = '  +16.15%  '

# Row 10: Polygon->BinanceUSD
$ws.Range("B10").Value = 'BinanceUSD'
$ws.Range("C10").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D10").Value = '''1.000'
$ws.Range("E10").Value = '  -0.93%  '

# Row 11: OKB->Polygon
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '''1.522'
$ws.Range("E11").Value = '  +9.57%  '

# Row 12: Dogecoin
$ws.Range("D12").Value = '''0.08774'
$ws.Range("E12").Value = '  +2.19%  '

# Row 13: Polkadot
$ws.Range("D13").Value = '''7.334'
$ws.Range("E13").Value = '  +14.16%  '

# Row 14: Solana
$ws.Range("D14").Value = '''23.08'
$ws.Range("E14").Value = '  +2.29%  '

# Row 15: ShibaInu
$ws.Range("D15").Value = '''0.00001318'
$ws.Range("E15").Value = '  +2.29%  '

# Row 16: Chainlink
$ws.Range("D16").Value = '''7.622'
$ws.Range("E16").Value = '  +7.65%  '

# Row 17: WrappedEther
$ws.Range("D17").Value = '1.692.99'
$ws.Range("E17").Value = '  +1.54%  '

# Row 18: Litecoin
$ws.Range("D18").Value = '''100.36'
$ws.Range("E18").Value = '  -0.40%  '

# Row 19: TRON
$ws.Range("D19").Value = '''0.07057'
$ws.Range("E19").Value = '  +4.59%  '

# Row 20: Avalanche
$ws.Range("D20").Value = '''19.46'
$ws.Range("E20").Value = '  +3.34%  '

# Row 21: Uniswap
$ws.Range("E21").Value = '  +1.77%  '

# Row 22: Dai
$ws.Range("E22").Value = '  -0.58%  '

# Row 23: Cosmos
$ws.Range("E23").Value = '  +4.33%  '

# Row 24: WrappedBTC
$ws.Range("D24").Value = '24.570.59'
$ws.Range("E24").Value = '  +1.97%  '

# Row 25: LidoDAOToken
$ws.Range("D25").Value = '''3.031'
$ws.Range("E25").Value = '  +13.24%  '

# Row 26: Toncoin
$ws.Range("D26").Value = '''2.310'
$ws.Range("E26").Value = '  +0.00%  '

# Row 27: EthereumClassic
$ws.Range("D27").Value = '''22.35'
$ws.Range("E27").Value = '  +3.49%  '

# Row 28: Monero
$ws.Range("D28").Value = '''159.73'
$ws.Range("E28").Value = '  +1.32%  '

# Row 29: HuobiToken
$ws.Range("D29").Value = '''5.177'
$ws.Range("E29").Value = '  -0.92%  '

# Row 30: BitcoinCash
$ws.Range("D30").Value = '''133.54'
$ws.Range("E30").Value = '  +2.44%  '

# Row 31: Filecoin
$ws.Range("D31").Value = '''7.627'
$ws.Range("E31").Value = '  +35.06%  '

# Row 32: WrappedliquidstakedEther2.0
$ws.Range("D32").Value = '1.880.89'
$ws.Range("E32").Value = '  +1.48%  '

# Row 33: ImmutableX
$ws.Range("D33").Value = '''1.090'
$ws.Range("E33").Value = '  -2.58%  '

# Row 34: Hedera->InternetComputer(DFINITY)
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '''7.375'
$ws.Range("E34").Value = '  +21.72%  '

# Row 35: InternetComputer(DFINITY)->Hedera
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.08586'
$ws.Range("E35").Value = '  +1.57%  '

# Row 36: WEMIXTOKEN
$ws.Range("E36").Value = '  +10.51%  '

# Row 37: FraxShare
$ws.Range("D37").Value = '''11.00'
$ws.Range("E37").Value = '  +7.14%  '

# Row 38: Algorand
$ws.Range("D38").Value = '''0.2715'
$ws.Range("E38").Value = '  +4.66%  '

# Row 39: Aptos
$ws.Range("D39").Value = '''14.72'
$ws.Range("E39").Value = '  +0.05%  '

# Row 40: VeChain
$ws.Range("D40").Value = '''0.02741'
$ws.Range("E40").Value = '  +10.87%  '

# Row 41: Stellar
$ws.Range("D41").Value = '''0.08992'
$ws.Range("E41").Value = '  +2.91%  '

# Row 42: TrustWalletToken
$ws.Range("D42").Value = '''1.468'
$ws.Range("E42").Value = '  +2.96%  '

# Row 43: TheSandbox
$ws.Range("D43").Value = '''0.7625'
$ws.Range("E43").Value = '  +4.97%  '

# Row 44: Decentraland
$ws.Range("D44").Value = '''0.7162'

# Row 45: EnergySwap
$ws.Range("D45").Value = '''15.39'
$ws.Range("E45").Value = '  +3.95%  '

# Row 46: NEARProtocol
$ws.Range("D46").Value = '''2.448'
$ws.Range("E46").Value = '  +4.80%  '

# Row 47: PancakeSwap
$ws.Range("D47").Value = '''4.173'
$ws.Range("E47").Value = '  +2.62%  '

# Row 48: Frax
$ws.Range("D48").Value = '''0.9996'
$ws.Range("E48").Value = '  -0.62%  '

# Row 49: Quant
$ws.Range("D49").Value = '''140.38'
$ws.Range("E49").Value = '  +1.44%  '

# Row 50: Flow
$ws.Range("D50").Value = '''1.324'
$ws.Range("E50").Value = '  +18.45%  '

# Row 51: BabyDogeCoin
$ws.Range("D51").Value = '''0.00000000384'
$ws.Range("E51").Value = '  +3.68%  '
